$d = $word.ActiveDocument

function Replace-ParagraphContent($Index, $InnerXml) {
    $p = $d.Paragraphs($Index)
    $r = $p.Range
    # Exclude the paragraph mark when this is the very last paragraph in the
    # document so we do not leave a stray empty trailing paragraph behind
    # (InsertXML always supplies a full <w:p>...</w:p>, and replacing a
    # range that runs right up to the final paragraph mark otherwise keeps
    # that old mark as an extra empty paragraph).
    if ($r.End -eq $d.Content.End) {
        $r = $d.Range($r.Start, $r.End - 1)
    }
    $xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xmlFrag)
}

function Append-Paragraphs($ParasXml) {
    $endRng = $d.Range($d.Content.End, $d.Content.End)
    $xmlFrag = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $ParasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $endRng.InsertXML($xmlFrag)
}

# 1) 'Matikka tyokalut' section paragraph: wrap 'sqrt' with proofErr spell markers
$p5Inner = '<w:r><w:t xml:space="preserve">Tämän jälkeen aloin työstämään laskinta omassa sivu XAML tiedostossa. Sain laskimen ilman liiallisia ongelmia valmiiksi vaikkakin edistyneemmät ominaisuudet veivät yllättävän kauan ongelmia, aiheutti </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sqrt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ominaisuudet ja toiminallisuus. Myöskään vielä ei ole varmuutta rajallisen testauksen vuoksi, että tuleeko ohjelma toimimaan pitkässä juoksussa. Implementoin näppäimistön niin kuin ajattelin suunnittelu vaiheessa mutta tein piilotus napin, joka on minusta hyvä ratkaisu antamaan enemmän vaihtoehtoa tarpeen mukaan. Loppu tulos on todella hyvä mutta varmasti hiomista vailla, johon palaan todennäköisesti myöhemmin.</w:t></w:r>'
Replace-ParagraphContent 5 $p5Inner

# 2) muunnin-sivu paragraph: wrap 'dictionary' with proofErr spell markers
$p6Inner = '<w:r><w:t xml:space="preserve">Seuraavaksi aloin kehittämään muunnin sivun. </w:t></w:r><w:r><w:t>Tämä alkoi samalla tavalla kuin laskin omassa XAML tiedostossa. Aluksi suunitelin rakenteen, En ole ihan täysin tyytyväinen lopputulokseen saatan palata siihen tyylittely vaiheessa</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Tämän jälkeen pureuduin muunnin </w:t></w:r><w:r><w:t xml:space="preserve">logiikkaan, jonka sain aika hyvin kehitettyä, vaikka lopputulos tuntui hieman liian epätehokaalle mutten usko, että on tehokkaampaa tapaa. Dictionary on minusta paras eteneminen tällaiseen ohjelmaan. Silti </w:t></w:r><w:r><w:t>tuntuu,</w:t></w:r><w:r><w:t xml:space="preserve"> että </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dictionary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> kohtia on liian monta, ehkä tulevaisuudessa palaan logiikkaan ja koitan parantaa </w:t></w:r><w:r><w:t>tätä,</w:t></w:r><w:r><w:t xml:space="preserve"> jos keksin uusia tapoja tuottaa tämän.</w:t></w:r>'
Replace-ParagraphContent 6 $p6Inner

# 3) Asetukset detail paragraph: wrap combboxit/label/checkboxit/boolean(x3)/app.xaml.cs/resourcedictionary/if
$p13Inner = '<w:r><w:t xml:space="preserve">Nyt kun kaikki alustus oli hoidettu aloin työstämään asetuksia. Aloitin miettimällä </w:t></w:r><w:r><w:t>hieman,</w:t></w:r><w:r><w:t xml:space="preserve"> miten tarkkaan ottaen haluan elementit </w:t></w:r><w:r><w:t>sijoittavat</w:t></w:r><w:r><w:t xml:space="preserve"> sivull</w:t></w:r><w:r><w:t xml:space="preserve">e. Kun sain </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>combboxit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>label</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ja </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checkboxit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> paikoilleen aloin implementoimaan logiikan. Opin että WPF nämä eivät ole kamalan </w:t></w:r><w:r><w:t>hankalia</w:t></w:r><w:r><w:t xml:space="preserve"> implementoida. </w:t></w:r><w:r><w:t xml:space="preserve">Tein erillisen asetukset tiedoston, joka on tietyn lainen tietokanta. Tähän tiedostoon tallensin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> muodossa Asetukset. Sitten vain yhdistin nämä elementit asetuksissa muutamaan näitä </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> muuttujia asetukset tiedostossa. Sitten tein tallenna nappi, joka kutsuu kaikkien asetuksien tallennus metodit. Tämän jälkeen minun piti implementoida</w:t></w:r><w:r><w:t xml:space="preserve"> näiden kyseisten muuttujien pohjalta asetus muutokset. Tämän implementoin </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>app.xaml.cs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> koodilla vaihtumaan </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>boolean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> mukaan kumman </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resourcedictionary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tiedoston tämä kutsuu. </w:t></w:r><w:r><w:t xml:space="preserve">Tein samalaiset </w:t></w:r><w:r><w:t>ratkaisut</w:t></w:r><w:r><w:t xml:space="preserve"> myös esimerkiksi ilmoituksille </w:t></w:r><w:r><w:t>yksinkertainen</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>toimi parhaiten tarkoitukseeni.</w:t></w:r>'
Replace-ParagraphContent 13 $p13Inner

# 4) Tyylittely paragraph: wrap resourcedictionary/xaml
$p15Inner = '<w:r><w:t xml:space="preserve">Viimeisimmiksi asioiksi jätin tyylitellyn. Rakenteet olivat jo aika hyvät omasta mielestäni mutta </w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve">loitin säätämällä niitä hieman. Seuraavaksi Tein </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>resourcedictionary</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xaml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>tiedostot,</w:t></w:r><w:r><w:t xml:space="preserve"> joissa määritelin normaalit tyylit ja värit elementeille. </w:t></w:r><w:r><w:t>Tyylittely</w:t></w:r><w:r><w:t xml:space="preserve"> oli paljon edestakaisin erilaisilla väreillä ja tyyleillä. Mutta loppuen lopuksi pääsin minusta</w:t></w:r><w:r><w:t xml:space="preserve"> hyviin tyyleihin, ainakin tumman teeman kanssa, vaaleasta teemasta en ole täysin varma mutta en kyllä yleensäkään pidä vaaleasta teemasta.</w:t></w:r>'
Replace-ParagraphContent 15 $p15Inner

# 5) Lisays paragraph: wrap checksum (x2) + rewrite final sentence
$p19Inner = '<w:r><w:t xml:space="preserve">Lisäsin seuraavaksi ohjelmaan vielä lisä ominaisuuksia, koska koin ettei ominaisuuksia välttämättä ole tarpeeksi. Täten lisäsin seuraavat ominaisuudet. Satunainen numero luoja, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checksum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> tiedosto tarkistus ja uuden työkalu elementti järjestelmä, jonka alle tein järjestelmän tieto tarkistus työkalun. Modulaarisen ja hyvän rakenteen vuoksi suurin osa näistä työkaluista meni todella helposti ja nopeasti. Vaikeuksia tuli järjestelmä tietojen kanssa. Loppuen lopuksi olen suurimalta osalta tyytyväinen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>checksum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> voisi tarvita pientä hiontaa. Eniten minusta järjestelmä tiedot </w:t></w:r><w:r><w:t>vaativat</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>hiomista ja parantamista,</w:t></w:r><w:r><w:t xml:space="preserve"> vaikka parhaillaan toimii ok tasolla.</w:t></w:r>'
Replace-ParagraphContent 19 $p19Inner

# 6) Append two new paragraphs: 'Lisays paranukset' heading + body paragraph
$p20 = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Lisäys paranukset</w:t></w:r></w:p>'
$p21 = '<w:p><w:r><w:t>Lisäyksien jälkeen aloin hiomaan näitä paranuksia. Aloitin saamalla järjestelmä tiedot</w:t></w:r><w:r><w:t xml:space="preserve"> päivittämään tiedot usein</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Seuraavaksi varmistin että päivittäminen toimii</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>erillisel</w:t></w:r><w:r><w:t xml:space="preserve">lä </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thread:illä</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> jottei ohjelma jäädy. Tämän jälkeen </w:t></w:r><w:r><w:t>varmistin,</w:t></w:r><w:r><w:t xml:space="preserve"> ettei järjestelmä tiedot päivitys </w:t></w:r><w:r><w:t>pyöri,</w:t></w:r><w:r><w:t xml:space="preserve"> kun tämä sivu ei ollut ladattuna</w:t></w:r><w:r><w:t>. Seuraavaksi paransin järjestelmä tiedot tyylittelyä, jotta tiedot olisivat helpompia lukea.</w:t></w:r></w:p>'
Append-Paragraphs ($p20 + $p21)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
